$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("A1:I1")
$rng.Interior.ThemeColor = 9
$rng.Interior.TintAndShade = 0.6
Write-Host "done"
